# This script applies the "Harvard case classification" recalculation update.
# The old "average_doctor" column (BP) is renamed to "average_doctor_old" and its
# historical values shift one column right into BQ; a freshly recomputed
# "average_doctor" series (and related recalculated precision/recall/F-score/etc.
# statistics across several model columns) replaces the values in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap the "average_doctor" / "average_doctor_old" labels ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Row 4: updated/recalculated statistics ---
$ws.Range("E4").Value = 0.474
$ws.Range("F4").Value = 0.05
$ws.Range("G4").Value = 0.224
$ws.Range("N4").Value = 0.461
$ws.Range("O4").Value = 0.062
$ws.Range("P4").Value = 0.248
$ws.Range("Q4").Value = 0.041
$ws.Range("R4").Value = 0.029
$ws.Range("S4").Value = 0.171
$ws.Range("W4").Value = 0.376
$ws.Range("X4").Value = 0.112
$ws.Range("Y4").Value = 0.335
$ws.Range("AI4").Value = 0.466
$ws.Range("AJ4").Value = 0.078
$ws.Range("AK4").Value = 0.28
$ws.Range("AU4").Value = 0.246
$ws.Range("BA4").Value = 2.049
$ws.Range("BB4").Value = 0.145
$ws.Range("BC4").Value = 0.38
$ws.Range("BG4").Value = 0.725
$ws.Range("BH4").Value = 0.14
$ws.Range("BI4").Value = 0.374
$ws.Range("BM4").Value = 0.743
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.6830000000000001
$ws.Range("BQ4").Value = 0.758

# --- Row 5: updated/recalculated statistics ---
$ws.Range("E5").Value = 0.609
$ws.Range("F5").Value = 0.05
$ws.Range("G5").Value = 0.224
$ws.Range("N5").Value = 0.717
$ws.Range("O5").Value = 0.077
$ws.Range("P5").Value = 0.278
$ws.Range("Q5").Value = 0.021
$ws.Range("R5").Value = 0.005
$ws.Range("S5").Value = 0.073
$ws.Range("W5").Value = 0.331
$ws.Range("X5").Value = 0.094
$ws.Range("Y5").Value = 0.307
$ws.Range("AI5").Value = 0.467
$ws.Range("AJ5").Value = 0.074
$ws.Range("AK5").Value = 0.272
$ws.Range("AU5").Value = 0.456
$ws.Range("AV5").Value = 0.081
$ws.Range("AW5").Value = 0.285
$ws.Range("BA5").Value = 1.291
$ws.Range("BB5").Value = 0.067
$ws.Range("BC5").Value = 0.259
$ws.Range("BG5").Value = 0.378
$ws.Range("BH5").Value = 0.049
$ws.Range("BI5").Value = 0.22
$ws.Range("BM5").Value = 0.522
$ws.Range("BN5").Value = 0.044
$ws.Range("BO5").Value = 0.209
$ws.Range("BP5").Value = 0.43
$ws.Range("BQ5").Value = 0.462

# --- Row 6: updated/recalculated statistics ---
$ws.Range("E6").Value = 0.533
$ws.Range("N6").Value = 0.5610000000000001
$ws.Range("Q6").Value = 0.028
$ws.Range("W6").Value = 0.352
$ws.Range("AI6").Value = 0.466
$ws.Range("AU6").Value = 0.32
$ws.Range("BA6").Value = 1.577
$ws.Range("BG6").Value = 0.497
$ws.Range("BM6").Value = 0.613
$ws.Range("BP6").Value = 0.526
$ws.Range("BQ6").Value = 0.571

# --- Row 7: updated/recalculated statistics ---
$ws.Range("E7").Value = 0.576
$ws.Range("N7").Value = 0.645
$ws.Range("Q7").Value = 0.023
$ws.Range("W7").Value = 0.339
$ws.Range("AI7").Value = 0.467
$ws.Range("AU7").Value = 0.39
$ws.Range("BA7").Value = 1.391
$ws.Range("BG7").Value = 0.418
$ws.Range("BM7").Value = 0.555
$ws.Range("BP7").Value = 0.464
$ws.Range("BQ7").Value = 0.5

# --- Row 8: updated/recalculated statistics ---
$ws.Range("E8").Value = 0.713
$ws.Range("F8").Value = 0.064
$ws.Range("G8").Value = 0.253
$ws.Range("N8").Value = 0.793
$ws.Range("O8").Value = 0.065
$ws.Range("P8").Value = 0.256
$ws.Range("Q8").Value = 0.022
$ws.Range("S8").Value = 0.111
$ws.Range("W8").Value = 0.418
$ws.Range("X8").Value = 0.119
$ws.Range("Y8").Value = 0.345
$ws.Range("AI8").Value = 0.548
$ws.Range("AJ8").Value = 0.117
$ws.Range("AK8").Value = 0.342
$ws.Range("AU8").Value = 0.416
$ws.Range("AV8").Value = 0.08599999999999999
$ws.Range("AW8").Value = 0.293
$ws.Range("BA8").Value = 1.768
$ws.Range("BB8").Value = 0.107
$ws.Range("BC8").Value = 0.326
$ws.Range("BG8").Value = 0.551
$ws.Range("BH8").Value = 0.111
$ws.Range("BI8").Value = 0.333
$ws.Range("BM8").Value = 0.6860000000000001
$ws.Range("BN8").Value = 0.06
$ws.Range("BO8").Value = 0.244
$ws.Range("BP8").Value = 0.589
$ws.Range("BQ8").Value = 0.624

# --- Row 9: updated/recalculated statistics ---
$ws.Range("E9").Value = 0.676
$ws.Range("F9").Value = 0.219
$ws.Range("G9").Value = 0.468
$ws.Range("N9").Value = 0.703
$ws.Range("O9").Value = 0.209
$ws.Range("P9").Value = 0.457
$ws.Range("W9").Value = 0.324
$ws.Range("X9").Value = 0.219
$ws.Range("Y9").Value = 0.468
$ws.Range("AI9").Value = 0.486
$ws.Range("AJ9").Value = 0.25
$ws.Range("AK9").Value = 0.5
$ws.Range("BA9").Value = 1.73
$ws.Range("BB9").Value = 0.25
$ws.Range("BC9").Value = 0.5
$ws.Range("BG9").Value = 0.5679999999999999
$ws.Range("BH9").Value = 0.245
$ws.Range("BI9").Value = 0.495
$ws.Range("BM9").Value = 0.676
$ws.Range("BN9").Value = 0.219
$ws.Range("BO9").Value = 0.468
$ws.Range("BP9").Value = 0.577
$ws.Range("BQ9").Value = 0.602

# --- Row 10: updated/recalculated statistics ---
$ws.Range("E10").Value = 0.8110000000000001
$ws.Range("F10").Value = 0.153
$ws.Range("G10").Value = 0.392
$ws.Range("N10").Value = 0.919
$ws.Range("O10").Value = 0.075
$ws.Range("P10").Value = 0.273
$ws.Range("W10").Value = 0.541
$ws.Range("AI10").Value = 0.595
$ws.Range("AJ10").Value = 0.241
$ws.Range("AK10").Value = 0.491
$ws.Range("AU10").Value = 0.432
$ws.Range("AV10").Value = 0.245
$ws.Range("AW10").Value = 0.495
$ws.Range("BA10").Value = 2.19
$ws.Range("BB10").Value = 0.209
$ws.Range("BC10").Value = 0.457
$ws.Range("BG10").Value = 0.649
$ws.Range("BH10").Value = 0.228
$ws.Range("BI10").Value = 0.477
$ws.Range("BM10").Value = 0.838
$ws.Range("BN10").Value = 0.136
$ws.Range("BO10").Value = 0.369
$ws.Range("BP10").Value = 0.73
$ws.Range("BQ10").Value = 0.754

# --- Row 11: updated/recalculated statistics ---
$ws.Range("E11").Value = 0.865
$ws.Range("F11").Value = 0.117
$ws.Range("G11").Value = 0.342
$ws.Range("N11").Value = 0.919
$ws.Range("O11").Value = 0.075
$ws.Range("P11").Value = 0.273
$ws.Range("W11").Value = 0.541
$ws.Range("AI11").Value = 0.676
$ws.Range("AJ11").Value = 0.219
$ws.Range("AK11").Value = 0.468
$ws.Range("AU11").Value = 0.5679999999999999
$ws.Range("AV11").Value = 0.245
$ws.Range("AW11").Value = 0.495
$ws.Range("BA11").Value = 2.19
$ws.Range("BB11").Value = 0.209
$ws.Range("BC11").Value = 0.457
$ws.Range("BG11").Value = 0.649
$ws.Range("BH11").Value = 0.228
$ws.Range("BI11").Value = 0.477
$ws.Range("BM11").Value = 0.838
$ws.Range("BN11").Value = 0.136
$ws.Range("BO11").Value = 0.369
$ws.Range("BP11").Value = 0.73
$ws.Range("BQ11").Value = 0.761

# --- Row 12: updated/recalculated statistics ---
$ws.Range("E12").Value = 1.438
$ws.Range("F12").Value = 0.9340000000000001
$ws.Range("G12").Value = 0.966
$ws.Range("N12").Value = 1.294
$ws.Range("O12").Value = 0.325
$ws.Range("P12").Value = 0.57
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 0.45
$ws.Range("Y12").Value = 0.671
$ws.Range("AI12").Value = 1.6
$ws.Range("AJ12").Value = 1.44
$ws.Range("AK12").Value = 1.2
$ws.Range("AU12").Value = 2.739
$ws.Range("AV12").Value = 3.497
$ws.Range("AW12").Value = 1.87
$ws.Range("BA12").Value = 3.816
$ws.Range("BB12").Value = 0.475
$ws.Range("BC12").Value = 0.6889999999999999
$ws.Range("BG12").Value = 1.167
$ws.Range("BH12").Value = 0.222
$ws.Range("BI12").Value = 0.471
$ws.Range("BM12").Value = 1.226
$ws.Range("BN12").Value = 0.239
$ws.Range("BO12").Value = 0.489
$ws.Range("BP12").Value = 1.272
$ws.Range("BQ12").Value = 1.271

# --- Row 13: updated/recalculated statistics ---
$ws.Range("E13").Value = 1.45
$ws.Range("F13").Value = 0.312
$ws.Range("G13").Value = 0.5590000000000001
$ws.Range("N13").Value = 1.769
$ws.Range("O13").Value = 0.458
$ws.Range("P13").Value = 0.677
$ws.Range("W13").Value = 0.971
$ws.Range("X13").Value = 0.202
$ws.Range("Y13").Value = 0.45
$ws.Range("AI13").Value = 1.154
$ws.Range("AJ13").Value = 0.303
$ws.Range("AK13").Value = 0.551
$ws.Range("AU13").Value = 2.056
$ws.Range("AV13").Value = 0.349
$ws.Range("AW13").Value = 0.591
$ws.Range("BA13").Value = 2.125
$ws.Range("BB13").Value = 0.267
$ws.Range("BC13").Value = 0.517
$ws.Range("BG13").Value = 0.538
$ws.Range("BH13").Value = 0.048
$ws.Range("BI13").Value = 0.219
$ws.Range("BM13").Value = 0.778
$ws.Range("BN13").Value = 0.13
$ws.Range("BO13").Value = 0.36
$ws.Range("BP13").Value = 0.708
$ws.Range("BQ13").Value = 0.668

Write-Host "Applied harvard case classification updates"
